$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.187.23"
$ws.Range("E2").Value = "'  +7.89%  "
$ws.Range("D3").Value = "'2.534.99"
$ws.Range("E3").Value = "'  +8.13%  "
$ws.Range("E4").Value = "'  +0.25%  "
$ws.Range("D5").Value = "'503.82"
$ws.Range("D6").Value = "'155.77"
$ws.Range("E6").Value = "'  +9.08%  "
$ws.Range("E7").Value = "'  +25.41%  "
$ws.Range("D8").Value = "'0.992"
$ws.Range("E8").Value = "'  -0.72%  "
$ws.Range("D9").Value = "'2.574.61"
$ws.Range("E9").Value = "'  +9.69%  "
$ws.Range("D10").Value = "'6.14"
$ws.Range("E10").Value = "'  +13.64%  "
$ws.Range("D11").Value = "'0.103"
$ws.Range("E11").Value = "'  +7.16%  "
$ws.Range("E12").Value = "'  +6.93%  "
$ws.Range("D13").Value = "'0.127"
$ws.Range("D14").Value = "'2.970.82"
$ws.Range("E14").Value = "'  +7.75%  "
$ws.Range("D15").Value = "'59.095.96"
$ws.Range("E15").Value = "'  +7.66%  "
$ws.Range("D16").Value = "'21.70"
$ws.Range("E16").Value = "'  +8.91%  "
$ws.Range("E17").Value = "'  +5.73%  "
$ws.Range("D18").Value = "'2.565.33"
$ws.Range("E18").Value = "'  +9.09%  "
$ws.Range("E19").Value = "'  +5.32%  "
$ws.Range("D20").Value = "'333.26"
$ws.Range("E20").Value = "'  +7.32%  "
$ws.Range("D21").Value = "'10.30"
$ws.Range("E21").Value = "'  +8.04%  "
$ws.Range("E22").Value = "'  +8.19%  "
$ws.Range("E23").Value = "'  +0.54%  "
$ws.Range("D24").Value = "'59.62"
$ws.Range("E24").Value = "'  +6.84%  "
$ws.Range("D25").Value = "'0.415"
$ws.Range("E25").Value = "'  +6.50%  "
$ws.Range("E26").Value = "'  +8.11%  "
$ws.Range("D27").Value = "'2.663.11"
$ws.Range("E27").Value = "'  +8.66%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "'  -0.52%  "
$ws.Range("D29").Value = "'0.0₃0825"
$ws.Range("E29").Value = "'  +10.45%  "
$ws.Range("D30").Value = "'7.39"
$ws.Range("E30").Value = "'  +3.82%  "
$ws.Range("E31").Value = "'  -0.26%  "
$ws.Range("D32").Value = "'157.41"
$ws.Range("E32").Value = "'  +7.15%  "
$ws.Range("D33").Value = "'19.27"
$ws.Range("E33").Value = "'  +7.48%  "
$ws.Range("E34").Value = "'  +7.53%  "
$ws.Range("E35").Value = "'  +9.79%  "
$ws.Range("E36").Value = "'  +9.96%  "
$ws.Range("D37").Value = "'3.88"
$ws.Range("E37").Value = "'  +9.81%  "
$ws.Range("D38").Value = "'0.846"
$ws.Range("E38").Value = "'  +3.55%  "
$ws.Range("D39").Value = "'3.74"
$ws.Range("E39").Value = "'  +12.21%  "
$ws.Range("E40").Value = "'  +8.21%  "
$ws.Range("D41").Value = "'35.09"
$ws.Range("E41").Value = "'  +5.10%  "
$ws.Range("D42").Value = "'290.08"
$ws.Range("E42").Value = "'  +15.37%  "
$ws.Range("D43").Value = "'0.102"
$ws.Range("E43").Value = "'  +7.82%  "
$ws.Range("D44").Value = "'0.623"
$ws.Range("E44").Value = "'  +8.51%  "
$ws.Range("D46").Value = "'0.989"
$ws.Range("E46").Value = "'  -0.98%  "
$ws.Range("D47").Value = "'0.761"
$ws.Range("E47").Value = "'  +21.85%  "
$ws.Range("B48").Value = "'RenderToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.83"
$ws.Range("E48").Value = "'  +10.42%  "
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'19.06"
$ws.Range("E49").Value = "'  +14.39%  "
$ws.Range("E50").Value = "'  +7.55%  "
$ws.Range("B51").Value = "'WhiteBITCoin"
$ws.Range("C51").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.24"
$ws.Range("E51").Value = "'  +0.88%  "
